$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures refreshed by the scraper run.
# Each cell holds its number as literal text (matching the original inline-string
# cells), so we force a temporary Text format while assigning the value and then
# clear the format again so no stray style survives on the cell.
$updates = @(
    @{ Cell = "D2"; Value = "303.53" }
    @{ Cell = "E2"; Value = "-2.08%" }
    @{ Cell = "D3"; Value = "35.57" }
    @{ Cell = "E3"; Value = "-0.22%" }
    @{ Cell = "D4"; Value = "5.079" }
    @{ Cell = "E4"; Value = "-0.83%" }
    @{ Cell = "D5"; Value = "0.08071" }
    @{ Cell = "E5"; Value = "-1.64%" }
    @{ Cell = "D6"; Value = "1.936" }
    @{ Cell = "E6"; Value = "-5.88%" }
    @{ Cell = "D7"; Value = "7.791" }
    @{ Cell = "E7"; Value = "-2.29%" }
    @{ Cell = "D8"; Value = "0.9271" }
    @{ Cell = "E8"; Value = "0.04%" }
    @{ Cell = "D9"; Value = "0.1541" }
    @{ Cell = "E9"; Value = "43.26%" }
    @{ Cell = "D10"; Value = "0.1897" }
    @{ Cell = "E10"; Value = "-1.36%" }
    @{ Cell = "D11"; Value = "0.08974" }
    @{ Cell = "E11"; Value = "-7.70%" }
    @{ Cell = "D12"; Value = "0.03456" }
    @{ Cell = "E12"; Value = "-4.02%" }
    @{ Cell = "D13"; Value = "0.09876" }
    @{ Cell = "E13"; Value = "-0.36%" }
    @{ Cell = "D14"; Value = "0.001422" }
    @{ Cell = "E14"; Value = "-1.34%" }
    @{ Cell = "D15"; Value = "0.005856" }
    @{ Cell = "E15"; Value = "-0.41%" }
    @{ Cell = "D16"; Value = "3.541" }
    @{ Cell = "E16"; Value = "1.91%" }
    @{ Cell = "D17"; Value = "4.053" }
    @{ Cell = "E17"; Value = "-1.90%" }
    @{ Cell = "E18"; Value = "3.26%" }
    @{ Cell = "D19"; Value = "0.3445" }
    @{ Cell = "E19"; Value = "0.56%" }
    @{ Cell = "D20"; Value = "0.1301" }
    @{ Cell = "E20"; Value = "-0.81%" }
    @{ Cell = "D21"; Value = "5.019" }
    @{ Cell = "E21"; Value = "-1.62%" }
    @{ Cell = "E22"; Value = "8.92%" }
    @{ Cell = "D23"; Value = "0.04491" }
    @{ Cell = "E23"; Value = "-1.35%" }
    @{ Cell = "D24"; Value = "0.001210" }
    @{ Cell = "E24"; Value = "-1.15%" }
    @{ Cell = "D25"; Value = "0.004812" }
    @{ Cell = "E25"; Value = "0.19%" }
    @{ Cell = "D26"; Value = "0.0001224" }
    @{ Cell = "E26"; Value = "-2.10%" }
    @{ Cell = "E27"; Value = "-32.30%" }
    @{ Cell = "D39"; Value = "0.01865" }
    @{ Cell = "E39"; Value = "-5.73%" }
    @{ Cell = "D40"; Value = "0.04803" }
    @{ Cell = "E40"; Value = "-2.30%" }
    @{ Cell = "D41"; Value = "0.01061" }
    @{ Cell = "E41"; Value = "7.93%" }
    @{ Cell = "D42"; Value = "0.007356" }
    @{ Cell = "E42"; Value = "-3.42%" }
    @{ Cell = "D43"; Value = "0.1346" }
    @{ Cell = "E43"; Value = "-2.75%" }
    @{ Cell = "D44"; Value = "0.002100" }
    @{ Cell = "E44"; Value = "-0.73%" }
    @{ Cell = "D45"; Value = "0.009714" }
    @{ Cell = "E45"; Value = "-15.91%" }
    @{ Cell = "D46"; Value = "0.00006226" }
    @{ Cell = "E46"; Value = "-4.26%" }
    @{ Cell = "D47"; Value = "0.00000000749" }
    @{ Cell = "E47"; Value = "-0.20%" }
    @{ Cell = "E48"; Value = "-63.13%" }
    @{ Cell = "D50"; Value = "0.00002097" }
    @{ Cell = "E50"; Value = "-0.20%" }
    @{ Cell = "D51"; Value = "0.0001997" }
    @{ Cell = "E51"; Value = "-0.20%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}

Write-Host "Updated $($updates.Count) cells"